$d = $word.ActiveDocument

foreach ($s in $d.Styles) {
    $f = $s.Font
    $f.LanguageID = "sv-SE"
    $f.LanguageIDFarEast = "en-US"
    $f.LanguageIDOther = "ar-SA"
}
